$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value while forcing text storage (matches source data
# which stores all these figures as literal strings, not numbers), without
# leaving a residual style change on the cell.
function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

Set-TextValue "D2" '71.117.71'
$ws.Range("E2").Value = '  +0.03%  '

Set-TextValue "D3" '3.836.86'
$ws.Range("E3").Value = '  +0.63%  '

Set-TextValue "D4" '0.999'
$ws.Range("E4").Value = '  -0.15%  '

Set-TextValue "D5" '711.48'
$ws.Range("E5").Value = '  +1.56%  '

Set-TextValue "D6" '173.06'
$ws.Range("E6").Value = '  -0.44%  '

Set-TextValue "D7" '3.835.46'
$ws.Range("E7").Value = '  +0.61%  '

$ws.Range("E8").Value = '  -0.05%  '

Set-TextValue "D10" '0.163'
$ws.Range("E10").Value = '  -0.02%  '

$ws.Range("E11").Value = '  +1.32%  '

$ws.Range("E12").Value = '  -0.09%  '

$ws.Range("E13").Value = '  -1.22%  '

Set-TextValue "D14" '36.74'
$ws.Range("E14").Value = '  +1.22%  '

Set-TextValue "D15" '4.477.88'
$ws.Range("E15").Value = '  +0.45%  '

Set-TextValue "D16" '3.788.49'
$ws.Range("E16").Value = '  -0.60%  '

Set-TextValue "D17" '71.006.71'
$ws.Range("E17").Value = '  -0.12%  '

Set-TextValue "D18" '7.22'
$ws.Range("E18").Value = '  +0.31%  '

$ws.Range("E19").Value = '  +0.67%  '

Set-TextValue "D20" '17.37'
$ws.Range("E20").Value = '  -2.58%  '

Set-TextValue "D21" '10.75'
$ws.Range("E21").Value = '  -3.96%  '

Set-TextValue "D22" '495.19'
$ws.Range("E22").Value = '  +3.16%  '

Set-TextValue "D23" '0.723'
$ws.Range("E23").Value = '  +1.47%  '

Set-TextValue "D24" '84.97'
$ws.Range("E24").Value = '  +1.24%  '

Set-TextValue "D26" '10.63'
$ws.Range("E26").Value = '  +1.75%  '

Set-TextValue "D27" '12.17'
$ws.Range("E27").Value = '  -1.38%  '

$ws.Range("B28").Value = 'PancakeSwap'
$ws.Range("C28").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue "D28" '3.20'
$ws.Range("E28").Value = '  +2.25%  '

$ws.Range("B29").Value = 'Fetch.AI'
$ws.Range("C29").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue "D29" '2.10'
$ws.Range("E29").Value = '  -3.12%  '

$ws.Range("E30").Value = '  -0.07%  '

Set-TextValue "D31" '7.51'
$ws.Range("E31").Value = '  -0.22%  '

Set-TextValue "D32" '2.26'
$ws.Range("E32").Value = '  -1.93%  '

Set-TextValue "D33" '29.48'
$ws.Range("E33").Value = '  -0.10%  '

Set-TextValue "D34" '0.177'
$ws.Range("E34").Value = '  -5.57%  '

Set-TextValue "D35" '9.20'
$ws.Range("E35").Value = '  -0.71%  '

Set-TextValue "D36" '3.797.55'
$ws.Range("E36").Value = '  +0.89%  '

$ws.Range("E37").Value = '  +0.02%  '

$ws.Range("E38").Value = '  -0.26%  '

Set-TextValue "D39" '2.31'
$ws.Range("E39").Value = '  +2.36%  '

Set-TextValue "D40" '6.02'
$ws.Range("E40").Value = '  +0.24%  '

$ws.Range("E41").Value = '  +5.28%  '

Set-TextValue "D42" '3.36'
$ws.Range("E42").Value = '  -2.40%  '

$ws.Range("E44").Value = '  +0.08%  '

Set-TextValue "D45" '0.000315'
$ws.Range("E45").Value = '  -3.62%  '

Set-TextValue "D46" '163.56'
$ws.Range("E46").Value = '  +0.22%  '

Set-TextValue "D47" '48.71'
$ws.Range("E47").Value = '  -0.58%  '

$ws.Range("E48").Value = '  +0.20%  '

Set-TextValue "D49" '415.47'
$ws.Range("E49").Value = '  +1.11%  '

Set-TextValue "D50" '8.62'
$ws.Range("E50").Value = '  +0.53%  '

$ws.Range("E51").Value = '  -1.55%  '
